# This script swaps the data of row 3 and row 4 for the columns that
# differ between the two rows (A, B, E, F, G, H, P, Q, R, S, AC), leaving
# all other columns (which are identical between the two rows) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New values for row 3 (previously found in row 4) ---
$ws.Range("A3").Value = 111634177
$ws.Range("B3").Value = 90350
$ws.Range("E3").Value = 4786
$ws.Range("F3").Value = "Mandelriska"
$ws.Range("G3").Value = "Lactarius volemus"
$ws.Range("H3").Value = "(Fr.:Fr.) Fr."
$ws.Range("P3").Value = "Hermansmåla, sydväst om Siggagölen, Bl"
$ws.Range("Q3").Value = 522930.7548289222
$ws.Range("R3").Value = 6247121.901725554
$ws.Range("S3").Value = 25
$ws.Range("AC3").Value = "Rikligt"

# --- New values for row 4 (previously found in row 3) ---
$ws.Range("A4").Value = 111634171
$ws.Range("B4").Value = 73683
$ws.Range("E4").Value = 306
$ws.Range("F4").Value = "Kornig nållav"
$ws.Range("G4").Value = "Chaenotheca chlorella"
$ws.Range("H4").Value = "(Ach.) Müll.Arg."
$ws.Range("P4").Value = "Hermansmåla, söder om Siggagölen, Bl"
$ws.Range("Q4").Value = 522996.846862453
$ws.Range("R4").Value = 6247111.736777187
$ws.Range("S4").Value = 10
$ws.Range("AC4").Value = "På askhögstubbe."
